$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value  = "0.4.0-snapshot-1"                 # Version
$meta.Range("B6").Value  = "draft"                             # Status
$meta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"         # Date
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"      # Contact

# --- Elements sheet: swap columns AK (37) and AL (38) ---
$elem = $wb.Worksheets.Item("Elements")

function Swap-CellValues($ws, $cellA, $cellB) {
    $valA = $ws.Range($cellA).Value2
    $valB = $ws.Range($cellB).Value2
    $ws.Range($cellA).Value = $valB
    $ws.Range($cellB).Value = $valA
}

# header row text ("Mapping: RIM Mapping" <-> "Mapping: Spécification métier ...")
Swap-CellValues $elem "AK1" "AL1"
# data rows (rows 2 and 4 hold identical/empty values in both columns, so no edit needed there)
Swap-CellValues $elem "AK3" "AL3"
Swap-CellValues $elem "AK5" "AL5"
Swap-CellValues $elem "AK6" "AL6"

# swap the column widths to match (AK becomes the wide column, AL the narrow one)
$elem.Columns.Item(37).ColumnWidth = 91.5
$elem.Columns.Item(38).ColumnWidth = 24.166666666666668
